# Update the "想去人数" (F column) figures that changed between the two
# data-refresh snapshots, on both the "展览" sheet and the aggregate
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1851
$ws1.Range("F7").Value  = 2502
$ws1.Range("F8").Value  = 157
$ws1.Range("F10").Value = 172
$ws1.Range("F11").Value = 1507
$ws1.Range("F12").Value = 521
$ws1.Range("F18").Value = 204
$ws1.Range("F22").Value = 153
$ws1.Range("F23").Value = 45
$ws1.Range("F24").Value = 1581
$ws1.Range("F27").Value = 571
$ws1.Range("F28").Value = 199
$ws1.Range("F30").Value = 406

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1851
$ws4.Range("F8").Value  = 2502
$ws4.Range("F9").Value  = 157
$ws4.Range("F11").Value = 172
$ws4.Range("F12").Value = 1507
$ws4.Range("F13").Value = 521
$ws4.Range("F19").Value = 204
$ws4.Range("F23").Value = 153
$ws4.Range("F24").Value = 45
$ws4.Range("F25").Value = 1581
$ws4.Range("F28").Value = 571
$ws4.Range("F29").Value = 199
$ws4.Range("F31").Value = 406

$wb.Save()
